# Agregue una pequeña parte del seaborn para el proyecto
#
# Changes applied:
#   1. Slide 3 ("CASO ELEGIDO"): reposition/resize the body text box and
#      switch it from non-wrapping to word-wrapped text (so the paragraph
#      text reflows instead of running off in one line).
#   2. Slide 10 ("REFERENCIAS"): move the "Librería pandas versión 2.1.4"
#      caption up, and remove the pandas logo picture that used to sit
#      above it.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 - "CASO ELEGIDO": resize/reposition the paragraph text box and turn
# on word-wrap (bodyPr wrap="none" -> wrap="square").
# ---------------------------------------------------------------------------
$s3  = $p.Slides.Item(3)
$tb3 = $s3.Shapes.Item("CuadroTexto 2")

$tb3.Left   = 46.507478714960634
$tb3.Top    = 82.58070756141731
$tb3.Width  = 866.9848328496063
$tb3.Height = 467.7234649669291

$tb3.TextFrame.WordWrap = [Microsoft.Office.Core.MsoTriState]::msoTrue

# ---------------------------------------------------------------------------
# Slide 10 - "REFERENCIAS": move the caption text box up and delete the
# pandas logo picture that used to sit above it.
# ---------------------------------------------------------------------------
$s10  = $p.Slides.Item(10)
$tb10 = $s10.Shapes.Item("CuadroTexto 2")

$tb10.Left = 44.03724489448819
$tb10.Top  = 118.23173148346457

$pic10 = $s10.Shapes.Item("Picture 2")
$pic10.Delete()
